# Append a new row (row 38) of data to each of the four worksheets,
# matching the structure/style of the existing rows (e.g. row 37).
#
# Column layout for every sheet:
#   A: timestamp (date-formatted number, same style as existing column A)
#   B: hex-byte string
#   C: hex-byte string
#   D: hex-byte string
#   E: hex-byte string
#   F: number
#   G: number
#   H: number
#   I: number

$wb = $excel.ActiveWorkbook

$rows = @(
    @{
        Sheet = "ROW35-FE-LIFTER"
        A = 45743.31901806713
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x7a"
        E = "0xd"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 378
        I = 13
    },
    @{
        Sheet = "ROW35-MID-LIFTER"
        A = 45743.16461385417
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x7e"
        E = "0xe"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 382
        I = 14
    },
    @{
        Sheet = "ROW02-FE-LIFTER"
        A = 45743.31371341435
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x7a"
        E = "0x3"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 378
        I = 3
    },
    @{
        Sheet = "ROW02-MID-LIFTER"
        A = 45743.36913008102
        B = "0x01,0x90"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x01,0x7e"
        E = "0x3"
        F = 400
        G = [double]"9.85046333984776e+23"
        H = 382
        I = 3
    }
)

foreach ($rowData in $rows) {
    $ws = $wb.Worksheets.Item($rowData.Sheet)

    # Find the first empty row in column A (the new row goes right after
    # the current last used row).
    $newRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

    $ws.Cells.Item($newRow, 1).Value = $rowData.A
    $ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($newRow, 2).Value = $rowData.B
    $ws.Cells.Item($newRow, 3).Value = $rowData.C
    $ws.Cells.Item($newRow, 4).Value = $rowData.D
    $ws.Cells.Item($newRow, 5).Value = $rowData.E
    $ws.Cells.Item($newRow, 6).Value = $rowData.F
    $ws.Cells.Item($newRow, 7).Value = $rowData.G
    $ws.Cells.Item($newRow, 8).Value = $rowData.H
    $ws.Cells.Item($newRow, 9).Value = $rowData.I
}

Write-Output "Appended row 38 to all four sheets"
